{"js": "// Indonesian Welcome.docx translation update:\n// Remove the trailing clause \", until Smartcash reaches a considerable\n// market cap\" from the SmartCash mining paragraph, leaving the sentence\n// ending at \"...for quite some time.\"\n\nconst body = context.document.body;\n\nconst results = body.search(\"ASICs have yet to be created for the Keccak mining algorithm\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const para = results.items[0].paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  const oldText = para.text;\n  const needle = \", until Smartcash reaches a considerable market cap\";\n  if (oldText.indexOf(needle) !== -1) {\n    const newText = oldText.replace(needle, \"\");\n    const range = para.getRange();\n    range.insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Indonesian Welcome.docx translation update:\n# Remove the trailing clause \", until Smartcash reaches a considerable\n# market cap\" from the SmartCash mining paragraph, leaving the sentence\n# ending at \"...for quite some time.\"\n\n$d = $word.ActiveDocument\n\n$needle = \", until Smartcash reaches a considerable market cap\"\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text.IndexOf(\"Keccak\") -ge 0 -and $r.Text.IndexOf($needle) -ge 0) {\n        $oldText = $r.Text\n        $newText = $oldText.Replace($needle, \"\")\n        # Paragraph.Range.Text includes the trailing paragraph mark (CR);\n        # strip it from the replacement text so we don't insert an extra\n        # paragraph break when writing back.\n        $newText = $newText.TrimEnd([char]13)\n        $r.Text = $newText\n        break\n    }\n}\n"}
